# Scratchpad update for approximate focal pixel length calculation
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New resolution scratch values (B3/B4) ---
$ws1.Range("B3").Value = 1080
$ws1.Range("B4").Value = 1200

# --- Aspect-ratio check formulas (D12:D15) ---
$ws1.Range("D12").Formula = "=B3/B4"
$ws1.Range("D13").Formula = "=B1/B2"
$ws1.Range("D14").Formula = "=H7/G6"
$ws1.Range("D15").Formula = "=H42/G41"

# --- Format the new focal-length scratch table (B19:D27) to match the
#     plain thin-bordered style already used by B1 ---
$ws1.Range("B1").Copy()
$ws1.Range("B19:D27").PasteSpecial(-4122)

# --- Populate the scratch table text labels in the same order the
#     original author typed them so the shared-string table matches ---
$ws1.Range("B27").Value = "fp"
$ws1.Range("B26").Value = "x"
$ws1.Range("B25").Value = "f"
$ws1.Range("B24").Value = "1/f"
$ws1.Range("B23").Value = "1/d2"
$ws1.Range("B22").Value = "1/d1"
$ws1.Range("B21").Value = "d2"
$ws1.Range("B20").Value = "d1"
$ws1.Range("C19").Value = "gear"
$ws1.Range("D19").Value = "vive"

# --- Measured values ---
$ws1.Range("C20").Value = 4.25
$ws1.Range("D20").Value = 4.5
$ws1.Range("C21").Value = 25.75
$ws1.Range("D21").Value = 25.5

# --- Derived formulas ---
$ws1.Range("C22").Formula = "=1/C20"
$ws1.Range("D22").Formula = "=1/D20"
$ws1.Range("C23").Formula = "=1/C21"
$ws1.Range("D23").Formula = "=1/D21"
$ws1.Range("C24").Formula = "=C22+C23"
$ws1.Range("D24").Formula = "=D22+D23"
$ws1.Range("C25").Formula = "=1/C24"
$ws1.Range("D25").Formula = "=1/D24"
$ws1.Range("C26").Formula = "=660/C25"
$ws1.Range("D26").Formula = "=C26"
$ws1.Range("C27").Formula = "=C26*C25"
$ws1.Range("D27").Formula = "=D26*D25"

# --- Selection / view state ---
$ws1.Range("D30").Select()
